$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shift D:M -> F:... wait, shift existing D:K to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formatting from column F (shifted original column D) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Set values for the two new quarter columns (D = FY2018-12-31 qtr, E = FY2018-09-30 qtr) ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1057700
$ws.Range("E8").Value = 1150400
$ws.Range("D9").Value = 676000
$ws.Range("E9").Value = 726500
$ws.Range("D10").Value = 381700
$ws.Range("E10").Value = 423900
$ws.Range("D12").Value = 43300
$ws.Range("E12").Value = 44800
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 53500
$ws.Range("E14").Value = 9700
$ws.Range("D15").Value = 65100
$ws.Range("E15").Value = 65800
$ws.Range("D17").Value = 1008400
$ws.Range("E17").Value = 1018200
$ws.Range("D18").Value = 49300
$ws.Range("E18").Value = 132200
$ws.Range("D20").Value = -38100
$ws.Range("E20").Value = -500
$ws.Range("D21").Value = 96000
$ws.Range("E21").Value = 226000
$ws.Range("D22").Value = 55400
$ws.Range("E22").Value = 66100
$ws.Range("D23").Value = -44200
$ws.Range("E23").Value = 65600
$ws.Range("D24").Value = -13100
$ws.Range("E24").Value = 1800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -31100
$ws.Range("E26").Value = 63800
$ws.Range("D27").Value = -31100
$ws.Range("E27").Value = 63800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 7800
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 38100
$ws.Range("E32").Value = 500
$ws.Range("D33").Value = -23300
$ws.Range("E33").Value = 63800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -23300
$ws.Range("E35").Value = 63800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 458200
$ws.Range("E41").Value = 352400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 810400
$ws.Range("E43").Value = 901100
$ws.Range("D44").Value = 473300
$ws.Range("E44").Value = 490800
$ws.Range("D45").Value = 135900
$ws.Range("E45").Value = 123300
$ws.Range("D46").Value = 1877800
$ws.Range("E46").Value = 1867500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 450900
$ws.Range("E48").Value = 445700
$ws.Range("D49").Value = 4204300
$ws.Range("E49").Value = 4279300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 97600
$ws.Range("E52").Value = 125700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 6630500
$ws.Range("E54").Value = 6718300
$ws.Range("D57").Value = 399200
$ws.Range("E57").Value = 441400
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 291400
$ws.Range("E59").Value = 323200
$ws.Range("D60").Value = 690600
$ws.Range("E60").Value = 764600
$ws.Range("D61").Value = 3985900
$ws.Range("E61").Value = 3983800
$ws.Range("D62").Value = 197200
$ws.Range("E62").Value = 214800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 4873800
$ws.Range("E66").Value = 4963200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -249800
$ws.Range("E72").Value = -226500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1756800
$ws.Range("E76").Value = 1755100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -23300
$ws.Range("E81").Value = 63800
$ws.Range("D83").Value = 84800
$ws.Range("E83").Value = 94300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 132300
$ws.Range("E89").Value = 226800
$ws.Range("D91").Value = -26900
$ws.Range("E91").Value = -24600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -22900
$ws.Range("E94").Value = -18100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 1000
$ws.Range("E100").Value = -400000
$ws.Range("D101").Value = -4600
$ws.Range("E101").Value = -2000
$ws.Range("D102").Value = 105800
$ws.Range("E102").Value = -193300

# --- Apply data corrections to shifted historical columns (source data revision) ---
$ws.Range("H9").Value = 1412600
$ws.Range("I9").Value = 1399300
$ws.Range("H10").Value = -292200
$ws.Range("I10").Value = -270500
$ws.Range("H12").Value = 90000
$ws.Range("I12").Value = 89100
$ws.Range("H17").Value = 1030200
$ws.Range("I17").Value = 1003300
$ws.Range("H18").Value = 90200
$ws.Range("I18").Value = 125500
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 4400
$ws.Range("H24").Value = -5000
$ws.Range("H26").Value = 31200
$ws.Range("H27").Value = 31200
$ws.Range("H29").Value = 22400
$ws.Range("H32").Value = -200
$ws.Range("I32").Value = -4400
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
